$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data: replace the "xxx" hotel name with real hotel names per city,
# and fix "London " (trailing space) to "London".
$ws.Range("B2").Value = "Jumeirah Beach Hotel"
$ws.Range("B3").Value = "Grand Plaza Apartments"
$ws.Range("A3").Value = "London"

# Column B needs to be widened to fit the longer hotel names.
$ws.Columns.Item(2).ColumnWidth = 21.26953125

# Move the active selection to A3.
$ws.Range("A3").Select()
